$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @{ Addr = "H12"; Val = 399.2 },
    @{ Addr = "I12"; Val = 399 },
    @{ Addr = "J12"; Val = 400 },
    @{ Addr = "K12"; Val = 399 },
    @{ Addr = "L12"; Val = 400 },
    @{ Addr = "M12"; Val = -229 },
    @{ Addr = "N12"; Val = -740 },
    @{ Addr = "H15"; Val = 2574.6414 },
    @{ Addr = "I15"; Val = 2574.6414 },
    @{ Addr = "K15"; Val = 7723.924199999999 },
    @{ Addr = "M15"; Val = -7554.924199999999 },
    @{ Addr = "H64"; Val = 3737.5 },
    @{ Addr = "J64"; Val = 3960 },
    @{ Addr = "L64"; Val = 3960 },
    @{ Addr = "N64"; Val = -4456 },
    @{ Addr = "H67"; Val = 3737.5 },
    @{ Addr = "J67"; Val = 3960 },
    @{ Addr = "L67"; Val = 3960 },
    @{ Addr = "N67"; Val = -5676 },
    @{ Addr = "H74"; Val = 3900 },
    @{ Addr = "I74"; Val = 4000 },
    @{ Addr = "K74"; Val = 4000 },
    @{ Addr = "M74"; Val = -3064 },
    @{ Addr = "H77"; Val = 3900 },
    @{ Addr = "I77"; Val = 4000 },
    @{ Addr = "K77"; Val = 20000 },
    @{ Addr = "M77"; Val = -15320 },
    @{ Addr = "H97"; Val = 452.5 },
    @{ Addr = "J97"; Val = 436.66666 },
    @{ Addr = "L97"; Val = 1309.99998 },
    @{ Addr = "N97"; Val = -2301.99998 },
    @{ Addr = "H107"; Val = 2830.0417 },
    @{ Addr = "I107"; Val = 2874.5334 },
    @{ Addr = "J107"; Val = 2755.889 },
    @{ Addr = "K107"; Val = 2874.5334 },
    @{ Addr = "L107"; Val = 2755.889 },
    @{ Addr = "M107"; Val = -954.5333999999998 },
    @{ Addr = "N107"; Val = -6595.889 },
    @{ Addr = "H120"; Val = 30000 },
    @{ Addr = "J120"; Val = 30000 },
    @{ Addr = "L120"; Val = 30000 },
    @{ Addr = "N120"; Val = -39676 },
    @{ Addr = "H137"; Val = 1398.9796 },
    @{ Addr = "J137"; Val = 1604.1666 },
    @{ Addr = "L137"; Val = 4812.4998 },
    @{ Addr = "N137"; Val = -9912.4998 },
    @{ Addr = "H138"; Val = 441057.1 },
    @{ Addr = "J138"; Val = 519567.78 },
    @{ Addr = "L138"; Val = 1558703.34 },
    @{ Addr = "N138"; Val = -1568983.34 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @{ Addr = "H2"; Val = 823.82355 },
    @{ Addr = "I2"; Val = 733.6667 },
    @{ Addr = "J2"; Val = 1500 },
    @{ Addr = "K2"; Val = 733.6667 },
    @{ Addr = "L2"; Val = 1500 },
    @{ Addr = "M2"; Val = -620.6667 },
    @{ Addr = "N2"; Val = -1726 },
    @{ Addr = "H32"; Val = 4294.6265 },
    @{ Addr = "I32"; Val = 4165.726 },
    @{ Addr = "K32"; Val = 4165.726 },
    @{ Addr = "M32"; Val = -3878.726 },
    @{ Addr = "H48"; Val = 249900 },
    @{ Addr = "J48"; Val = 249900 },
    @{ Addr = "L48"; Val = 249900 },
    @{ Addr = "N48"; Val = -250668 },
    @{ Addr = "H108"; Val = 29700 },
    @{ Addr = "J108"; Val = 29700 },
    @{ Addr = "L108"; Val = 29700 },
    @{ Addr = "N108"; Val = -37380 },
    @{ Addr = "H116"; Val = 823.82355 },
    @{ Addr = "I116"; Val = 733.6667 },
    @{ Addr = "J116"; Val = 1500 },
    @{ Addr = "K116"; Val = 733.6667 },
    @{ Addr = "L116"; Val = 1500 },
    @{ Addr = "M116"; Val = 1560.3333 },
    @{ Addr = "N116"; Val = -6088 },
    @{ Addr = "H122"; Val = 2961.3572 },
    @{ Addr = "I122"; Val = 2593.7273 },
    @{ Addr = "J122"; Val = 4309.3335 },
    @{ Addr = "K122"; Val = 7781.1819 },
    @{ Addr = "L122"; Val = 12928.0005 },
    @{ Addr = "M122"; Val = -5331.1819 },
    @{ Addr = "N122"; Val = -17828.0005 },
    @{ Addr = "H132"; Val = 3444.25 },
    @{ Addr = "I132"; Val = 2990.25 },
    @{ Addr = "J132"; Val = 4125.25 },
    @{ Addr = "K132"; Val = 8970.75 },
    @{ Addr = "L132"; Val = 12375.75 },
    @{ Addr = "M132"; Val = -6440.75 },
    @{ Addr = "N132"; Val = -17435.75 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @{ Addr = "H3"; Val = 823.82355 },
    @{ Addr = "I3"; Val = 733.6667 },
    @{ Addr = "J3"; Val = 1500 },
    @{ Addr = "K3"; Val = 733.6667 },
    @{ Addr = "L3"; Val = 1500 },
    @{ Addr = "M3"; Val = -619.6667 },
    @{ Addr = "N3"; Val = -1728 },
    @{ Addr = "H86"; Val = 2314.6956 },
    @{ Addr = "I86"; Val = 2566.25 },
    @{ Addr = "J86"; Val = 2040.2727 },
    @{ Addr = "K86"; Val = 2566.25 },
    @{ Addr = "L86"; Val = 2040.2727 },
    @{ Addr = "M86"; Val = -1443.25 },
    @{ Addr = "N86"; Val = -4286.2727 },
    @{ Addr = "H89"; Val = 2314.6956 },
    @{ Addr = "I89"; Val = 2566.25 },
    @{ Addr = "J89"; Val = 2040.2727 },
    @{ Addr = "K89"; Val = 12831.25 },
    @{ Addr = "L89"; Val = 10201.3635 },
    @{ Addr = "M89"; Val = -7215.25 },
    @{ Addr = "N89"; Val = -21433.3635 },
    @{ Addr = "H94"; Val = 15625320 },
    @{ Addr = "I94"; Val = 27778054 },
    @{ Addr = "J94"; Val = 377.14285 },
    @{ Addr = "K94"; Val = 27778054 },
    @{ Addr = "L94"; Val = 377.14285 },
    @{ Addr = "M94"; Val = -27777603 },
    @{ Addr = "N94"; Val = -1279.14285 },
    @{ Addr = "H134"; Val = 1273.7778 },
    @{ Addr = "I134"; Val = 1002.4167 },
    @{ Addr = "J134"; Val = 1816.5 },
    @{ Addr = "K134"; Val = 3007.2501 },
    @{ Addr = "L134"; Val = 5449.5 },
    @{ Addr = "M134"; Val = -472.2501000000002 },
    @{ Addr = "N134"; Val = -10519.5 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @{ Addr = "H31"; Val = 1353.1428 },
    @{ Addr = "I31"; Val = 1211.3334 },
    @{ Addr = "K31"; Val = 1211.3334 },
    @{ Addr = "M31"; Val = -916.3334 },
    @{ Addr = "H34"; Val = 1353.1428 },
    @{ Addr = "I34"; Val = 1211.3334 },
    @{ Addr = "K34"; Val = 1211.3334 },
    @{ Addr = "M34"; Val = -1009.3334 },
    @{ Addr = "H86"; Val = 7469819 },
    @{ Addr = "I86"; Val = 13383733 },
    @{ Addr = "J86"; Val = 77427 },
    @{ Addr = "K86"; Val = 13383733 },
    @{ Addr = "L86"; Val = 77427 },
    @{ Addr = "M86"; Val = -13382610 },
    @{ Addr = "N86"; Val = -79673 },
    @{ Addr = "H89"; Val = 7469819 },
    @{ Addr = "I89"; Val = 13383733 },
    @{ Addr = "J89"; Val = 77427 },
    @{ Addr = "K89"; Val = 66918665 },
    @{ Addr = "L89"; Val = 387135 },
    @{ Addr = "M89"; Val = -66913049 },
    @{ Addr = "N89"; Val = -398367 },
    @{ Addr = "H122"; Val = 1723.625 },
    @{ Addr = "I122"; Val = 1684.1428 },
    @{ Addr = "K122"; Val = 5052.428400000001 },
    @{ Addr = "M122"; Val = -2602.428400000001 },
    @{ Addr = "H132"; Val = 1666.2 },
    @{ Addr = "I132"; Val = 1440.7368 },
    @{ Addr = "J132"; Val = 2055.6365 },
    @{ Addr = "K132"; Val = 4322.2104 },
    @{ Addr = "L132"; Val = 6166.9095 },
    @{ Addr = "M132"; Val = -1792.2104 },
    @{ Addr = "N132"; Val = -11226.9095 },
    @{ Addr = "H141"; Val = 682955.5600000001 },
    @{ Addr = "J141"; Val = 682955.5600000001 },
    @{ Addr = "L141"; Val = 682955.5600000001 },
    @{ Addr = "N141"; Val = -693315.5600000001 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @{ Addr = "H5"; Val = 957.0606 },
    @{ Addr = "I5"; Val = 1060.6923 },
    @{ Addr = "J5"; Val = 572.1429000000001 },
    @{ Addr = "K5"; Val = 3182.0769 },
    @{ Addr = "L5"; Val = 1716.4287 },
    @{ Addr = "M5"; Val = -3070.0769 },
    @{ Addr = "N5"; Val = -1940.4287 },
    @{ Addr = "H131"; Val = 13158788 },
    @{ Addr = "J131"; Val = 941.89856 },
    @{ Addr = "L131"; Val = 2825.69568 },
    @{ Addr = "N131"; Val = -12905.69568 },
    @{ Addr = "H135"; Val = 957.0606 },
    @{ Addr = "I135"; Val = 1060.6923 },
    @{ Addr = "J135"; Val = 572.1429000000001 },
    @{ Addr = "K135"; Val = 9546.2307 },
    @{ Addr = "L135"; Val = 5149.2861 },
    @{ Addr = "M135"; Val = -7011.2307 },
    @{ Addr = "N135"; Val = -10219.2861 },
    @{ Addr = "H140"; Val = 23537.646 },
    @{ Addr = "I140"; Val = 41913.36 },
    @{ Addr = "J140"; Val = 3564.0435 },
    @{ Addr = "K140"; Val = 125740.08 },
    @{ Addr = "L140"; Val = 10692.1305 },
    @{ Addr = "M140"; Val = -120560.08 },
    @{ Addr = "N140"; Val = -21052.1305 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @{ Addr = "H97"; Val = 530.05 },
    @{ Addr = "I97"; Val = 503.85715 },
    @{ Addr = "J97"; Val = 591.1667 },
    @{ Addr = "K97"; Val = 503.85715 },
    @{ Addr = "L97"; Val = 591.1667 },
    @{ Addr = "M97"; Val = -7.85714999999999 },
    @{ Addr = "N97"; Val = -1583.1667 },
    @{ Addr = "H122"; Val = 2876.182 },
    @{ Addr = "I122"; Val = 3103.0667 },
    @{ Addr = "K122"; Val = 9309.2001 },
    @{ Addr = "M122"; Val = -6859.2001 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @{ Addr = "H61"; Val = 1148.0667 },
    @{ Addr = "I61"; Val = 1015.63635 },
    @{ Addr = "J61"; Val = 1512.25 },
    @{ Addr = "K61"; Val = 1015.63635 },
    @{ Addr = "L61"; Val = 1512.25 },
    @{ Addr = "M61"; Val = -813.63635 },
    @{ Addr = "N61"; Val = -1916.25 },
    @{ Addr = "H102"; Val = 16500 },
    @{ Addr = "J102"; Val = 16500 },
    @{ Addr = "L102"; Val = 16500 },
    @{ Addr = "N102"; Val = -22990 },
    @{ Addr = "H113"; Val = 1148.0667 },
    @{ Addr = "I113"; Val = 1015.63635 },
    @{ Addr = "J113"; Val = 1512.25 },
    @{ Addr = "K113"; Val = 1015.63635 },
    @{ Addr = "L113"; Val = 1512.25 },
    @{ Addr = "M113"; Val = 1154.36365 },
    @{ Addr = "N113"; Val = -5852.25 },
    @{ Addr = "H132"; Val = 2940.1365 },
    @{ Addr = "I132"; Val = 2648.9285 },
    @{ Addr = "J132"; Val = 3449.75 },
    @{ Addr = "K132"; Val = 7946.7855 },
    @{ Addr = "L132"; Val = 10349.25 },
    @{ Addr = "M132"; Val = -5416.7855 },
    @{ Addr = "N132"; Val = -15409.25 },
    @{ Addr = "H136"; Val = 1137.5 },
    @{ Addr = "I136"; Val = 1137.5 },
    @{ Addr = "J136"; Val = 0 },
    @{ Addr = "K136"; Val = 3412.5 },
    @{ Addr = "L136"; Val = 0 },
    @{ Addr = "M136"; Val = -862.5 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }
foreach ($a in @("N136")) { $ws.Range($a).ClearContents() }

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @{ Addr = "H16"; Val = 47000 },
    @{ Addr = "J16"; Val = 47000 },
    @{ Addr = "L16"; Val = 47000 },
    @{ Addr = "N16"; Val = -47584 },
    @{ Addr = "H106"; Val = 0 },
    @{ Addr = "J106"; Val = 0 },
    @{ Addr = "L106"; Val = 0 },
    @{ Addr = "H107"; Val = 448.875 },
    @{ Addr = "I107"; Val = 391.2 },
    @{ Addr = "K107"; Val = 1173.6 },
    @{ Addr = "M107"; Val = 746.4000000000001 },
    @{ Addr = "H113"; Val = 435.7 },
    @{ Addr = "I113"; Val = 257.125 },
    @{ Addr = "K113"; Val = 771.375 },
    @{ Addr = "M113"; Val = 1398.625 },
    @{ Addr = "H121"; Val = 0 },
    @{ Addr = "J121"; Val = 0 },
    @{ Addr = "L121"; Val = 0 },
    @{ Addr = "H126"; Val = 76926380 },
    @{ Addr = "I126"; Val = 83336090 },
    @{ Addr = "K126"; Val = 250008270 },
    @{ Addr = "M126"; Val = -250005800 }
)
foreach ($e in $edits) { $ws.Range($e.Addr).Value = $e.Val }
foreach ($a in @("N106","N121")) { $ws.Range($a).ClearContents() }
